$d = $word.ActiveDocument

$table = $d.Tables.Item(1)
$row = $table.Rows.Item(4)

$row.Cells.Item(1).Range.Text = "1.2"
$row.Cells.Item(1).Range.Font.Size = 12

$row.Cells.Item(2).Range.Text = "19/03/2021"
$row.Cells.Item(2).Range.Font.Size = 12

$row.Cells.Item(3).Range.Text = "Abhinav Jain"
$row.Cells.Item(3).Range.Font.Size = 12

$row.Cells.Item(4).Range.Text = "No change. For git demo."
$row.Cells.Item(4).Range.Font.Size = 12
